$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 1.02
    "C2" = 1.035139555456449
    "D2" = 1.044580604364526
    "E2" = 1.034206980081355
    "F2" = 1.055109735890809
    "I2" = 1.039309114280873
    "J2" = 1.040254962063045
    "K2" = 1.04735125367426
    "L2" = 1.037007118210116
    "M2" = 1.057851097053617
    "N2" = 1.01728512677642
    "B3" = 1.02
    "C3" = 1.036370131472121
    "D3" = 1.045557119118181
    "E3" = 1.035261092962709
    "F3" = 1.0562312821505
    "I3" = 1.03962524413423
    "J3" = 1.041127542890279
    "K3" = 1.048138604227128
    "L3" = 1.037869718098207
    "M3" = 1.058785233101543
    "N3" = 1.017581686209386
    "B4" = 1.02
    "C4" = 1.037166021957374
    "D4" = 1.046188371854336
    "E4" = 1.035943168612179
    "F4" = 1.056956648957779
    "I4" = 1.039827900192209
    "J4" = 1.041691300690423
    "K4" = 1.048646841340944
    "L4" = 1.038427277039346
    "M4" = 1.059388746635347
    "N4" = 1.017773103647507
    "B5" = 1.02
    "C5" = 1.037500526940488
    "D5" = 1.046453603599611
    "E5" = 1.036229912989397
    "F5" = 1.057261510876728
    "I5" = 1.03991264206607
    "J5" = 1.041928099738918
    "K5" = 1.048860210380778
    "L5" = 1.038661532087075
    "M5" = 1.059642240933073
    "N5" = 1.017853461836531
    "B6" = 1.02
    "C6" = 1.037556686741268
    "D6" = 1.046498128589631
    "E6" = 1.036278058653024
    "F6" = 1.05731269374298
    "I6" = 1.039926843945142
    "J6" = 1.041967847371469
    "K6" = 1.048896018777614
    "L6" = 1.038700856202932
    "M6" = 1.059684790706471
    "N6" = 1.017866947672264
    "B7" = 1.02
    "C7" = 1.037170491969827
    "D7" = 1.046191916469965
    "E7" = 1.035947000104005
    "F7" = 1.056960722859545
    "I7" = 1.039829034303589
    "J7" = 1.04169446561382
    "K7" = 1.048649693541942
    "L7" = 1.038430407726211
    "M7" = 1.059392134712705
    "N7" = 1.017774177843915
    "B8" = 1.02
    "C8" = 1.035555513292803
    "D8" = 1.044910750552875
    "E8" = 1.034563224392717
    "F8" = 1.055488840593178
    "I8" = 1.039416345718166
    "J8" = 1.040550033725438
    "K8" = 1.047617598007954
    "L8" = 1.037298762968958
    "M8" = 1.058166986920434
    "N8" = 1.017385449258664
    "B9" = 1.02
    "C9" = 1.032706744124424
    "D9" = 1.042648398916056
    "E9" = 1.032124736374104
    "F9" = 1.052892472589973
    "I9" = 1.038674557405657
    "J9" = 1.038526753326217
    "K9" = 1.04578944020567
    "L9" = 1.035300007799772
    "M9" = 1.056000912052863
    "N9" = 1.016696796314646
    "B10" = 1.02
    "C10" = 1.030805422045235
    "D10" = 1.041136890974214
    "E10" = 1.030498934019153
    "F10" = 1.051159651495331
    "I10" = 1.038170202935204
    "J10" = 1.037173351056334
    "K10" = 1.044564235885082
    "L10" = 1.033964303083939
    "M10" = 1.054551946667423
    "N10" = 1.016235208868961
    "B11" = 1.02
    "C11" = 1.029981587245481
    "D11" = 1.040481599665299
    "E11" = 1.029794891465436
    "F11" = 1.050408849608929
    "I11" = 1.037949473912958
    "J11" = 1.036586215274754
    "K11" = 1.044032168310862
    "L11" = 1.033385151855425
    "M11" = 1.053923346385824
    "N11" = 1.016034741600638
    "B12" = 1.02
    "C12" = 1.029675492764137
    "D12" = 1.040238073915556
    "E12" = 1.029533368129133
    "F12" = 1.050129894700071
    "I12" = 1.037867133075933
    "J12" = 1.036367959372184
    "K12" = 1.043834301058768
    "L12" = 1.0331699102126
    "M12" = 1.053689676132999
    "N12" = 1.015960188963921
    "B13" = 1.02
    "C13" = 1.029741154988587
    "D13" = 1.040290316531385
    "E13" = 1.029589466278245
    "F13" = 1.050189734835715
    "I13" = 1.037884811415725
    "J13" = 1.036414783612987
    "K13" = 1.043876754852322
    "L13" = 1.033216085689885
    "M13" = 1.053739807370935
    "N13" = 1.015976184845838
    "B14" = 1.02
    "C14" = 1.029956287142604
    "D14" = 1.040461472235307
    "E14" = 1.029773274090424
    "F14" = 1.050385792621676
    "I14" = 1.037942674783323
    "J14" = 1.036568177601313
    "K14" = 1.044015817316864
    "L14" = 1.033367362352547
    "M14" = 1.053904034813112
    "N14" = 1.0160285808987
    "B15" = 1.02
    "C15" = 1.030088825737431
    "D15" = 1.040566910786204
    "E15" = 1.0298865226976
    "F15" = 1.05050658048255
    "I15" = 1.037978279617395
    "J15" = 1.036662666449129
    "K15" = 1.044101467286212
    "L15" = 1.033460553069734
    "M15" = 1.054005196846836
    "N15" = 1.0160608518732
    "B16" = 1.02
    "C16" = 1.030860085412159
    "D16" = 1.041180363557111
    "E16" = 1.030545657611716
    "F16" = 1.051209469508187
    "I16" = 1.038184802631745
    "J16" = 1.037212293925088
    "K16" = 1.044599514766187
    "L16" = 1.034002722828857
    "M16" = 1.054593639622848
    "N16" = 1.016248500603927
    "B17" = 1.02
    "C17" = 1.031343726951055
    "D17" = 1.041564951513098
    "E17" = 1.030959098455713
    "F17" = 1.05165024380249
    "I17" = 1.03831372205781
    "J17" = 1.037556764026669
    "K17" = 1.044911511898722
    "L17" = 1.034342601057621
    "M17" = 1.054962434902089
    "N17" = 1.01636604754956
    "B18" = 1.02
    "C18" = 1.031625774032461
    "D18" = 1.041789198226042
    "E18" = 1.031200245887152
    "F18" = 1.051907293899747
    "I18" = 1.038388692827321
    "J18" = 1.03775758123857
    "K18" = 1.045093345474934
    "L18" = 1.034540770907229
    "M18" = 1.055177432361031
    "N18" = 1.016434553140103
    "B19" = 1.02
    "C19" = 1.031721935923698
    "D19" = 1.041865647553199
    "E19" = 1.03128247002615
    "F19" = 1.05199493357318
    "I19" = 1.038414217650723
    "J19" = 1.03782603672739
    "K19" = 1.04515532078168
    "L19" = 1.034608328967293
    "M19" = 1.055250721548034
    "N19" = 1.016457902028773
    "B20" = 1.02
    "C20" = 1.031291842279773
    "D20" = 1.04152369683418
    "E20" = 1.030914740759994
    "F20" = 1.051602957673008
    "I20" = 1.038299913579164
    "J20" = 1.037519816660526
    "K20" = 1.044878052989321
    "L20" = 1.034306143148393
    "M20" = 1.054922878523257
    "N20" = 1.016353441831746
    "B21" = 1.02
    "C21" = 1.029892938470814
    "D21" = 1.040411074508702
    "E21" = 1.029719147576581
    "F21" = 1.050328060537437
    "I21" = 1.037925645192938
    "J21" = 1.036523011533514
    "K21" = 1.043974873348628
    "L21" = 1.033322818445917
    "M21" = 1.053855678909187
    "N21" = 1.016013154060132
    "B22" = 1.02
    "C22" = 1.029012895654942
    "D22" = 1.039710821703763
    "E22" = 1.028967368622809
    "F22" = 1.049526055086346
    "I22" = 1.037688289461032
    "J22" = 1.035895310446622
    "K22" = 1.043405655851826
    "L22" = 1.032703873380533
    "M22" = 1.053183645270952
    "N22" = 1.015798679703257
    "B23" = 1.02
    "C23" = 1.029479471052198
    "D23" = 1.040082106081364
    "E23" = 1.029365907310358
    "F23" = 1.049951254400667
    "I23" = 1.037814309659155
    "J23" = 1.036228159203651
    "K23" = 1.043707537496805
    "L23" = 1.033032053779448
    "M23" = 1.053540002385444
    "N23" = 1.015912426257183
    "B24" = 1.02
    "C24" = 1.031315286892625
    "D24" = 1.041542338285598
    "E24" = 1.030934784110763
    "F24" = 1.051624324382013
    "I24" = 1.038306153733853
    "J24" = 1.037536511913601
    "K24" = 1.044893172090919
    "L24" = 1.034322617141186
    "M24" = 1.054940752702292
    "N24" = 1.016359137990878
    "B25" = 1.02
    "C25" = 1.033443586193899
    "D25" = 1.043233843285439
    "E25" = 1.03275516358749
    "F25" = 1.053564026678776
    "I25" = 1.038868057853607
    "J25" = 1.039050615584411
    "K25" = 1.046263191929202
    "L25" = 1.035817292773379
    "M25" = 1.056561755591769
    "N25" = 1.016875266043119
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Output ("Updated {0} cells" -f $updates.Count)
